$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ส้ม"
$ws.Range("B3").Value = 5

$ws.Range("A4").Value = "โรตี"
$ws.Range("B4").Value = 7

$ws.Range("A5").Value = "พะแนง"
$ws.Range("B5").Value = 45

$ws.Range("A6").Value = "หมูปิ้ง"
$ws.Range("B6").Value = 9

$ws.Range("A7").Value = "แหนม"
$ws.Range("B7").Value = 15

$ws.Range("C5").Select()
